$d = $word.ActiveDocument

$d.Content.Find.Execute("295÷4=73, 3", $true, $false, $false, $false, $false, $true, 1, $false, "320÷6=53, 2", 2)
$d.Content.Find.Execute("234÷7=33, 3", $true, $false, $false, $false, $false, $true, 1, $false, "417÷9=46, 3", 2)
$d.Content.Find.Execute("939÷9=104, 3", $true, $false, $false, $false, $false, $true, 1, $false, "723÷8=90, 3", 2)
$d.Content.Find.Execute("204÷3=68, 0", $true, $false, $false, $false, $false, $true, 1, $false, "310÷3=103, 1", 2)
$d.Content.Find.Execute("587÷4=146, 3", $true, $false, $false, $false, $false, $true, 1, $false, "947÷4=236, 3", 2)
$d.Content.Find.Execute("291÷4=72, 3", $true, $false, $false, $false, $false, $true, 1, $false, "273÷4=68, 1", 2)
$d.Content.Find.Execute("517÷9=57, 4", $true, $false, $false, $false, $false, $true, 1, $false, "885÷8=110, 5", 2)
$d.Content.Find.Execute("317÷4=79, 1", $true, $false, $false, $false, $false, $true, 1, $false, "127÷4=31, 3", 2)
$d.Content.Find.Execute("130÷5=26, 0", $true, $false, $false, $false, $false, $true, 1, $false, "874÷8=109, 2", 2)
$d.Content.Find.Execute("745÷8=93, 1", $true, $false, $false, $false, $false, $true, 1, $false, "666÷5=133, 1", 2)
$d.Content.Find.Execute("454÷5=90, 4", $true, $false, $false, $false, $false, $true, 1, $false, "419÷8=52, 3", 2)
$d.Content.Find.Execute("406÷3=135, 1", $true, $false, $false, $false, $false, $true, 1, $false, "732÷8=91, 4", 2)
$d.Content.Find.Execute("163÷8=20, 3", $true, $false, $false, $false, $false, $true, 1, $false, "383÷4=95, 3", 2)
$d.Content.Find.Execute("851÷5=170, 1", $true, $false, $false, $false, $false, $true, 1, $false, "762÷4=190, 2", 2)
$d.Content.Find.Execute("894÷9=99, 3", $true, $false, $false, $false, $false, $true, 1, $false, "479÷5=95, 4", 2)
$d.Content.Find.Execute("445÷4=111, 1", $true, $false, $false, $false, $false, $true, 1, $false, "789÷8=98, 5", 2)
$d.Content.Find.Execute("501÷3=167, 0", $true, $false, $false, $false, $false, $true, 1, $false, "316÷4=79, 0", 2)
$d.Content.Find.Execute("473÷7=67, 4", $true, $false, $false, $false, $false, $true, 1, $false, "475÷6=79, 1", 2)
$d.Content.Find.Execute("868÷4=217, 0", $true, $false, $false, $false, $false, $true, 1, $false, "929÷8=116, 1", 2)
$d.Content.Find.Execute("623÷8=77, 7", $true, $false, $false, $false, $false, $true, 1, $false, "562÷2=281, 0", 2)
$d.Content.Find.Execute("172÷4=43, 0", $true, $false, $false, $false, $false, $true, 1, $false, "270÷6=45, 0", 2)
$d.Content.Find.Execute("482÷2=241, 0", $true, $false, $false, $false, $false, $true, 1, $false, "480÷8=60, 0", 2)
$d.Content.Find.Execute("963÷6=160, 3", $true, $false, $false, $false, $false, $true, 1, $false, "950÷7=135, 5", 2)
$d.Content.Find.Execute("896÷8=112, 0", $true, $false, $false, $false, $false, $true, 1, $false, "229÷2=114, 1", 2)
$d.Content.Find.Execute("855÷4=213, 3", $true, $false, $false, $false, $false, $true, 1, $false, "598÷5=119, 3", 2)
